$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell K1 with same style as other header cells (A1:J1)
$ws.Range("K1").Value = "intervention_type"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats

# Add new data values in column K for rows 2-4
$ws.Range("K2").Value = "DEVICE"
$ws.Range("K3").Value = "PROCEDURE"
$ws.Range("K4").Value = "PROCEDURE"
